# Fruta / hortaliza, semanal
#
# A new weekly price-report row for "Femacal de La Calera" / "Chirimoya" is
# inserted at row 92, pushing the existing rows 92-172 down to 93-173.
# The worksheet's used range grows from A1:T172 to A1:T173.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92 (shifts rows 92..172 down to 93..173,
# carrying their existing content and formatting along with them).
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new record.
$ws.Range("A92").Value = 3
$ws.Range("B92").Value = "Femacal de La Calera"
$ws.Range("C92").Value = "Coquimbo"
$ws.Range("D92").Value = 44778
$ws.Range("E92").Value = 5
$ws.Range("F92").Value = "Fruta"
$ws.Range("G92").Value = 100107
$ws.Range("H92").Value = "Otros"
$ws.Range("I92").Value = 100107002
$ws.Range("J92").Value = "Chirimoya"
$ws.Range("K92").Value = "Cultivar IV Región"
$ws.Range("L92").Value = "Primera"
$ws.Range("M92").Value = 92
$ws.Range("N92").Value = 28000
$ws.Range("O92").Value = 30000
$ws.Range("P92").Value = 28978
$ws.Range("Q92").Value = "$/bandeja 10 kilos"
$ws.Range("R92").Value = "Provincia del Elquí"
$ws.Range("S92").Value = 2898
$ws.Range("T92").Value = 10
